$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so that numeric-looking
# strings (e.g. "259.38", "1.003") are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "26.284.39"
$ws.Range("E2").Value = "  -2.87%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.832.05"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "259.38"
$ws.Range("E5").Value = "  -8.10%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.5188"
$ws.Range("E7").Value = "  -1.95%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3212"
$ws.Range("E8").Value = "  -8.91%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06722"
$ws.Range("E9").Value = "  -4.05%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "18.75"
$ws.Range("E10").Value = "  -7.78%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.7668"
$ws.Range("E11").Value = "  -6.18%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07673"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.868.63"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "88.70"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.007"
$ws.Range("E15").Value = "  -3.71%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  -3.73%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007838"
$ws.Range("E19").Value = "  -4.12%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "26.334.19"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.091.04"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.526"
$ws.Range("E22").Value = "  -5.32%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "9.388"
$ws.Range("E23").Value = "  -7.87%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "5.868"
$ws.Range("E24").Value = "  -6.05%  "

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  -3.80%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "144.63"
$ws.Range("E26").Value = "  -1.09%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "1.664"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "16.90"
$ws.Range("E28").Value = "  -3.95%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "110.82"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "4.159"
$ws.Range("E30").Value = "  -5.47%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08717"
$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.087"
$ws.Range("E32").Value = "  -6.70%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.04819"
$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.125"
$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.857"
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6798"
$ws.Range("E36").Value = "  -8.95%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "3.086"
$ws.Range("E37").Value = "  -6.42%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01777"
$ws.Range("E38").Value = "  -5.63%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.204"
$ws.Range("E39").Value = "  -8.89%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.4898"
$ws.Range("E40").Value = "  -8.01%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "111.47"
$ws.Range("E41").Value = "  -4.96%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8976"
$ws.Range("E42").Value = "  -8.30%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.142"
$ws.Range("E43").Value = "  -3.07%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "7.688"
$ws.Range("E45").Value = "  -6.28%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4175"
$ws.Range("E46").Value = "  -9.42%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1252"
$ws.Range("E47").Value = "  -8.53%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05881"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.006"
$ws.Range("E49").Value = "  -4.81%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "35.21"
$ws.Range("E50").Value = "  -4.23%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "59.11"
$ws.Range("E51").Value = "  -4.45%  "
